# Atualizacao de bases das ligas, do dia: 06-04-2024 as 15:39
#
# The upstream scraper re-pulled results for "Paraguay Division Profesional".
# For three pairs of fixtures (rows 130/131, 134/135, 140/141) the two
# matches on the same date had been stored swapped (the stats for the first
# game were under the second game's row and vice versa); the row/id column
# (A) is correct and untouched, but every other column (B:AC) needs to be
# exchanged between the two rows. Four later fixtures (rows 210-213) simply
# got refreshed closing-odds figures (columns N:V) from the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rows 130 / 131 (id 128 / 129) : swap all match data ----
$ws.Range("B130").Value  = 7453204
$ws.Range("F130").Value  = "Cerro Porteno"
$ws.Range("G130").Value  = "Tacuary"
$ws.Range("H130").Value  = 1
$ws.Range("J130").Value  = "D"
$ws.Range("K130").Value  = 1.285
$ws.Range("L130").Value  = 5
$ws.Range("M130").Value  = 8
$ws.Range("N130").Value  = 1.285
$ws.Range("O130").Value  = 4.75
$ws.Range("P130").Value  = 8
$ws.Range("Q130").Value  = -1.5
$ws.Range("R130").Value  = 1.9
$ws.Range("S130").Value  = 1.9
$ws.Range("T130").Value  = 3
$ws.Range("U130").Value  = 1.9
$ws.Range("V130").Value  = 1.9
$ws.Range("X130").Value  = 3.75
$ws.Range("Y130").Value  = -1
$ws.Range("AA130").Value = 0.8999999999999999
$ws.Range("AC130").Value = 0.8999999999999999

$ws.Range("B131").Value  = 7454842
$ws.Range("F131").Value  = "Sportivo Luqueno"
$ws.Range("G131").Value  = "Libertad Asuncion"
$ws.Range("H131").Value  = 0
$ws.Range("J131").Value  = "A"
$ws.Range("K131").Value  = 4
$ws.Range("L131").Value  = 3.6
$ws.Range("M131").Value  = 1.727
$ws.Range("N131").Value  = 3.5
$ws.Range("O131").Value  = 3.3
$ws.Range("P131").Value  = 1.95
$ws.Range("Q131").Value  = 0.5
$ws.Range("R131").Value  = 1.8
$ws.Range("S131").Value  = 2
$ws.Range("T131").Value  = 2.5
$ws.Range("U131").Value  = 1.975
$ws.Range("V131").Value  = 1.825
$ws.Range("X131").Value  = -1
$ws.Range("Y131").Value  = 0.95
$ws.Range("AA131").Value = 1
$ws.Range("AC131").Value = 0.825

# ---- Rows 134 / 135 (id 132 / 133) : swap all match data ----
$ws.Range("B134").Value  = 7493427
$ws.Range("F134").Value  = "Tacuary"
$ws.Range("G134").Value  = "Sportivo Luqueno"
$ws.Range("H134").Value  = 1
$ws.Range("J134").Value  = "D"
$ws.Range("K134").Value  = 3.4
$ws.Range("L134").Value  = 3.3
$ws.Range("M134").Value  = 2
$ws.Range("N134").Value  = 3.2
$ws.Range("O134").Value  = 3.25
$ws.Range("P134").Value  = 2.1
$ws.Range("Q134").Value  = 0.25
$ws.Range("R134").Value  = 2.025
$ws.Range("S134").Value  = 1.775
$ws.Range("T134").Value  = 2.5
$ws.Range("U134").Value  = 1.975
$ws.Range("V134").Value  = 1.825
$ws.Range("W134").Value  = -1
$ws.Range("X134").Value  = 2.25
$ws.Range("Z134").Value  = 0.5125
$ws.Range("AA134").Value = -0.5
$ws.Range("AB134").Value = -1
$ws.Range("AC134").Value = 0.825

$ws.Range("B135").Value  = 7493428
$ws.Range("F135").Value  = "Guairena FC"
$ws.Range("G135").Value  = "Resistencia FC"
$ws.Range("H135").Value  = 4
$ws.Range("J135").Value  = "H"
$ws.Range("K135").Value  = 1.727
$ws.Range("L135").Value  = 3.6
$ws.Range("M135").Value  = 4.2
$ws.Range("N135").Value  = 1.45
$ws.Range("O135").Value  = 4.2
$ws.Range("P135").Value  = 6
$ws.Range("Q135").Value  = -1
$ws.Range("R135").Value  = 1.775
$ws.Range("S135").Value  = 2.025
$ws.Range("T135").Value  = 2.75
$ws.Range("U135").Value  = 1.825
$ws.Range("V135").Value  = 1.975
$ws.Range("W135").Value  = 0.45
$ws.Range("X135").Value  = -1
$ws.Range("Z135").Value  = 0.7749999999999999
$ws.Range("AA135").Value = -1
$ws.Range("AB135").Value = 0.825
$ws.Range("AC135").Value = -1

# ---- Rows 140 / 141 (id 138 / 139) : swap all match data ----
$ws.Range("B140").Value  = 7493310
$ws.Range("F140").Value  = "Libertad Asuncion"
$ws.Range("G140").Value  = "Tacuary"
$ws.Range("H140").Value  = 1
$ws.Range("J140").Value  = "A"
$ws.Range("K140").Value  = 1.363
$ws.Range("L140").Value  = 5
$ws.Range("M140").Value  = 7
$ws.Range("N140").Value  = 1.571
$ws.Range("O140").Value  = 4.2
$ws.Range("P140").Value  = 4.75
$ws.Range("Q140").Value  = -0.75
$ws.Range("R140").Value  = 1.8
$ws.Range("S140").Value  = 2
$ws.Range("T140").Value  = 2.75
$ws.Range("U140").Value  = 1.8
$ws.Range("V140").Value  = 2
$ws.Range("W140").Value  = -1
$ws.Range("Y140").Value  = 3.75
$ws.Range("Z140").Value  = -1
$ws.Range("AA140").Value = 1
$ws.Range("AB140").Value = 0.4
$ws.Range("AC140").Value = -0.5

$ws.Range("B141").Value  = 7493431
$ws.Range("F141").Value  = "Sportivo Trinidense"
$ws.Range("G141").Value  = "Guairena FC"
$ws.Range("H141").Value  = 7
$ws.Range("J141").Value  = "H"
$ws.Range("K141").Value  = 2.05
$ws.Range("L141").Value  = 3.3
$ws.Range("M141").Value  = 3.3
$ws.Range("N141").Value  = 2.6
$ws.Range("O141").Value  = 3.1
$ws.Range("P141").Value  = 2.6
$ws.Range("Q141").Value  = 0
$ws.Range("R141").Value  = 1.925
$ws.Range("S141").Value  = 1.875
$ws.Range("T141").Value  = 2.5
$ws.Range("U141").Value  = 2
$ws.Range("V141").Value  = 1.8
$ws.Range("W141").Value  = 1.6
$ws.Range("Y141").Value  = -1
$ws.Range("Z141").Value  = 0.925
$ws.Range("AA141").Value = -1
$ws.Range("AB141").Value = 1
$ws.Range("AC141").Value = -1

# ---- Rows 210-213 (id 208-211) : refreshed closing-odds values only ----
$ws.Range("N210").Value = 4.2
$ws.Range("O210").Value = 4
$ws.Range("P210").Value = 1.65
$ws.Range("Q210").Value = 0.75
$ws.Range("R210").Value = 1.95
$ws.Range("S210").Value = 1.85
$ws.Range("T210").Value = 2.5
$ws.Range("U210").Value = 1.825
$ws.Range("V210").Value = 1.975

$ws.Range("N211").Value = 1.7
$ws.Range("O211").Value = 3.6
$ws.Range("P211").Value = 4.75
$ws.Range("Q211").Value = -0.75
$ws.Range("R211").Value = 1.9
$ws.Range("S211").Value = 1.9
$ws.Range("U211").Value = 2
$ws.Range("V211").Value = 1.8

$ws.Range("R212").Value = 1.775
$ws.Range("S212").Value = 2.025

$ws.Range("R213").Value = 1.9
$ws.Range("S213").Value = 1.9
